$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.034808581472853
$ws.Range("D2").Value = 1.037965196692606
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.044784549568577
$ws.Range("I2").Value = 1.0398784209957
$ws.Range("J2").Value = 1.039925798267606
$ws.Range("K2").Value = 1.040754577511155
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.047554625497415
$ws.Range("N2").Value = 1.017243946977731

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.035782025253119
$ws.Range("D3").Value = 1.038696607897317
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.046006071013783
$ws.Range("I3").Value = 1.040180226626452
$ws.Range("J3").Value = 1.040542407058362
$ws.Range("K3").Value = 1.041296112483841
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.048586385825185
$ws.Range("N3").Value = 1.017450762914287

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.036411696154983
$ws.Range("D4").Value = 1.039169660032088
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.046796646435889
$ws.Range("I4").Value = 1.040374071264215
$ws.Range("J4").Value = 1.040940582141858
$ws.Range("K4").Value = 1.041645635323712
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.049253610684283
$ws.Range("N4").Value = 1.01758425108268

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.036676359035053
$ws.Range("D5").Value = 1.03936847797328
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.047129045819882
$ws.Range("I5").Value = 1.040455217799082
$ws.Range("J5").Value = 1.041107780402928
$ws.Range("K5").Value = 1.04179236225274
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.049534018615837
$ws.Range("N5").Value = 1.017640289072448

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.036720794121987
$ws.Range("D6").Value = 1.039401857246097
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.047184859598362
$ws.Range("I6").Value = 1.040468822381277
$ws.Range("J6").Value = 1.041135842314383
$ws.Range("K6").Value = 1.041816985883133
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.049581094884836
$ws.Range("N6").Value = 1.017649693377782

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.036415232792739
$ws.Range("D7").Value = 1.039172316856793
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.046801087808583
$ws.Range("I7").Value = 1.040375156906883
$ws.Range("J7").Value = 1.040942817017811
$ws.Range("K7").Value = 1.041647596731301
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.049257357874906
$ws.Range("N7").Value = 1.017585000181226

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.035137605355532
$ws.Range("D8").Value = 1.038212425720164
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.045197333951436
$ws.Range("I8").Value = 1.039980716504306
$ws.Range("J8").Value = 1.040134352148895
$ws.Range("K8").Value = 1.040937775200891
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.04790339553987
$ws.Range("N8").Value = 1.017313910814651

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.032884633068856
$ws.Range("D9").Value = 1.036519315083748
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.042372569848454
$ws.Range("I9").Value = 1.039274608766722
$ws.Range("J9").Value = 1.038703520240852
$ws.Range("K9").Value = 1.039680199718999
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.045514498663458
$ws.Range("N9").Value = 1.016833650862174

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.031381547949484
$ws.Range("D10").Value = 1.035389488607237
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.040490175938246
$ws.Range("I10").Value = 1.038796442931618
$ws.Range("J10").Value = 1.037745460039334
$ws.Range("K10").Value = 1.038837269885932
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.043919807981901
$ws.Range("N10").Value = 1.016511757367993

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.030730430944001
$ws.Range("D11").Value = 1.034900009294138
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.039675249465099
$ws.Range("I11").Value = 1.038587630576147
$ws.Range("J11").Value = 1.037329620152978
$ws.Range("K11").Value = 1.038471195493951
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.043228781519635
$ws.Range("N11").Value = 1.016371966718132

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.030488535988928
$ws.Range("D12").Value = 1.034718156644959
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.039372572531634
$ws.Range("I12").Value = 1.038509803396228
$ws.Range("J12").Value = 1.037175009373295
$ws.Range("K12").Value = 1.038335056935931
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.042972025075467
$ws.Range("N12").Value = 1.016319980883997

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.030540425143944
$ws.Range("D13").Value = 1.034757166379961
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.039437496740619
$ws.Range("I13").Value = 1.038526509581844
$ws.Range("J13").Value = 1.037208180669036
$ws.Range("K13").Value = 1.038364266459046
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.043027103771849
$ws.Range("N13").Value = 1.016331134798184

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.030710436672231
$ws.Range("D14").Value = 1.034884978075341
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.039650229633223
$ws.Range("I14").Value = 1.038581202762815
$ws.Range("J14").Value = 1.037316843024404
$ws.Range("K14").Value = 1.038459945547368
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.043207559558514
$ws.Range("N14").Value = 1.016367670804466

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.030815180927895
$ws.Range("D15").Value = 1.034963722023295
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.039781304401821
$ws.Range("I15").Value = 1.038614865923804
$ws.Range("J15").Value = 1.037383773698
$ws.Range("K15").Value = 1.038518875096282
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.043318733896588
$ws.Range("N15").Value = 1.016390173713559

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.031424754695624
$ws.Range("D16").Value = 1.03542196835451
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.040544263231203
$ws.Range("I16").Value = 1.03881026396355
$ws.Range("J16").Value = 1.037773037025391
$ws.Range("K16").Value = 1.038861542284075
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.043965658191194
$ws.Range("N16").Value = 1.016521026206101

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.031807051346795
$ws.Range("D17").Value = 1.035709345631411
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.041022889688896
$ws.Range("I17").Value = 1.038932359850504
$ws.Range("J17").Value = 1.038016945591449
$ws.Range("K17").Value = 1.03907619919804
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.044371318091743
$ws.Range("N17").Value = 1.01660299713484

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.032030012369936
$ws.Range("D18").Value = 1.035876943114233
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.041302080176134
$ws.Range("I18").Value = 1.039003406176762
$ws.Range("J18").Value = 1.038159117478171
$ws.Range("K18").Value = 1.039201300715613
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.044607883044136
$ws.Range("N18").Value = 1.016650769960601

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.032106031893145
$ws.Range("D19").Value = 1.035934085303337
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.041397279634309
$ws.Range("I19").Value = 1.039027602280821
$ws.Range("J19").Value = 1.038207578161171
$ws.Range("K19").Value = 1.039243939419328
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.044688537199013
$ws.Range("N19").Value = 1.016667052570618

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.031766037240946
$ws.Range("D20").Value = 1.035678515333857
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.040971535977952
$ws.Range("I20").Value = 1.038919277715852
$ws.Range("J20").Value = 1.037990786425584
$ws.Range("K20").Value = 1.039053179322995
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.044327799762055
$ws.Range("N20").Value = 1.016594206508135

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.030660373652257
$ws.Range("D21").Value = 1.034847341769273
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.039587584473499
$ws.Range("I21").Value = 1.038565104299139
$ws.Range("J21").Value = 1.037284848791909
$ws.Range("K21").Value = 1.038431774918451
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.043154422001429
$ws.Range("N21").Value = 1.016356913554913

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.029964961051932
$ws.Range("D22").Value = 1.034324528906155
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.038717571874832
$ws.Range("I22").Value = 1.03834088833649
$ws.Range("J22").Value = 1.036840133766285
$ws.Range("K22").Value = 1.038040134724564
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.042416218412947
$ws.Range("N22").Value = 1.016207363046864

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.030333635110774
$ws.Range("D23").Value = 1.034601702727272
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.039178769900423
$ws.Range("I23").Value = 1.038459894793582
$ws.Range("J23").Value = 1.037075967538174
$ws.Range("K23").Value = 1.038247839548741
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.042807597541963
$ws.Range("N23").Value = 1.016286676241438

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.031784569830527
$ws.Range("D24").Value = 1.035692446294886
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.040994740460276
$ws.Range("I24").Value = 1.038925189495002
$ws.Range("J24").Value = 1.038002606922677
$ws.Range("K24").Value = 1.039063581335605
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.044347463974923
$ws.Range("N24").Value = 1.016598178735755

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.033467274328299
$ws.Range("D25").Value = 1.036957218360921
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.04310269669122
$ws.Range("I25").Value = 1.039458463972274
$ws.Range("J25").Value = 1.039074160355678
$ws.Range("K25").Value = 1.040006115271765
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.046132451157413
$ws.Range("N25").Value = 1.016958113124724

Write-Host "Applied all vm_pu.xlsx changes (380 kV case)"
